$wb = $excel.ActiveWorkbook

# --- Rename sheets (tab names) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778195466595"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778212398808"
$wb.Worksheets.Item(3).Name = "RS_TO-1650477821241881"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778212888823"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778213679123"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778195086646.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778195296962.csv"
$ws1.Range("B4").Value = "go_stims-16504778195316596.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778195456955.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_5-16504778197766976.csv"
$ws2.Range("B3").Value = "TB-16504778212268782.csv"
$ws2.Range("B4").Value = "OB-16504778209119103.csv"
$ws2.Range("B5").Value = "OB-1650477820257876.csv"
$ws2.Range("B6").Value = "TB-16504778211318781.csv"
$ws2.Range("B7").Value = "TB-16504778209969108.csv"
$ws2.Range("B8").Value = "OB-16504778202096612.csv"
$ws2.Range("B9").Value = "ZB-match_4-16504778200026975.csv"
$ws2.Range("B10").Value = "ZB-match_7-16504778197316608.csv"

# --- Sheet 3 (RS) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778212558823.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778212438774.csv"
$ws4.Range("B4").Value = "MM_stims-16504778212719107.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778212568784.csv"
$ws4.Range("B6").Value = "MM_stims-16504778212879114.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778212728767.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16504778213528774.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778213358777.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778212958808.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778213199117.csv"
